$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to already-populated rows (just the "new cases" count) ---
$ws.Range("C595").Value = 54
$ws.Range("C596").Value = 28
$ws.Range("C597").Value = 29

# --- Newly populated rows 598-601 ---
# Row 598
$ws.Range("C598").Value = 35
$ws.Range("E598").Value = 2
$ws.Range("F598").Value = 0
$ws.Range("G598").Value = 6

# Row 599
$ws.Range("C599").Value = 22
$ws.Range("E599").Value = 2
$ws.Range("F599").Value = 1
$ws.Range("G599").Value = 6

# Row 600
$ws.Range("C600").Value = 21
$ws.Range("E600").Value = 1
$ws.Range("F600").Value = 1
$ws.Range("G600").Value = 8

# Row 601
$ws.Range("C601").Value = 5
$ws.Range("E601").Value = 2
$ws.Range("F601").Value = 1
$ws.Range("G601").Value = 9

# --- L/M columns are formatted as Text (numFmt "@") in this sheet, but the
# historical data stores genuine numbers there. Writing straight into a
# Text-formatted cell makes Excel coerce the input to a text string, so we
# flip the number format to a plain numeric one, write the value, then
# restore the original "@" (Text) display format -- this keeps the cell's
# style index the same as before while the stored value stays numeric.
foreach ($r in 598..601) {
    $cellL = $ws.Range("L$r")
    $cellL.NumberFormat = "0"
    $cellL.Value = 0
    $cellL.NumberFormat = "@"

    $cellM = $ws.Range("M$r")
    $cellM.NumberFormat = "0"
    $cellM.Value = 0
    $cellM.NumberFormat = "@"
}
